# Commit: change "la roi" to "le roi" in the French TechReport colophon.
#
# The colophon paragraph reads (before):
#   "© Sa Majesté la Roi du chef du Canada, représenté par le ministre
#    du ministère des Pêches et des Océans, ..."
# and should read (after):
#   "© Sa Majesté le Roi du chef du Canada, représenté par le ministre
#    du ministère des Pêches et des Océans, ..."
#
# i.e. the feminine article "la" before "Roi" is corrected to the
# masculine article "le".

$d = $word.ActiveDocument

$d.Content.Find.Execute("la Roi", $true, $false, $false, $false, $false,
                         $true, 1, $false, "le Roi", 2) | Out-Null
